# Update the cryptocurrency price snapshot (commit: "Updated symbol list on
# Wed Jan 11 00:27:47 UTC 2023 with GitHub Actions").
# Columns: D=Price, E=Volume(1h), F=Data, G=Hora.
# All four columns are stored as plain text in the sheet, so each cell's
# NumberFormat is forced to "@" (Text) before the value is assigned —
# otherwise Excel would auto-coerce strings like "277.67", "1.20%" or
# "11-1-2023" into a number/percentage/date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colMap = @{ D = 4; E = 5; F = 6; G = 7 }

$updates = @(
    @{ Row=2; D="277.67"; E="1.20%"; F="11-1-2023"; G="0" },
    @{ Row=3; D="27.17"; E="1.90%"; F="11-1-2023"; G="0" },
    @{ Row=4; D="4.877"; E="-0.43%"; F="11-1-2023"; G="0" },
    @{ Row=5; D="0.06417"; E="1.44%"; F="11-1-2023"; G="0" },
    @{ Row=6; D="6.957"; E="0.67%"; F="11-1-2023"; G="0" },
    @{ Row=7; D="1.254"; E="-4.79%"; F="11-1-2023"; G="0" },
    @{ Row=8; D="0.8833"; E="0.13%"; F="11-1-2023"; G="0" },
    @{ Row=9; D="0.1518"; E="3.80%"; F="11-1-2023"; G="0" },
    @{ Row=10; D="0.05112"; E="0.51%"; F="11-1-2023"; G="0" },
    @{ Row=11; D="0.07523"; E="1.97%"; F="11-1-2023"; G="0" },
    @{ Row=12; D="0.02966"; E="-6.12%"; F="11-1-2023"; G="0" },
    @{ Row=13; D="0.09005"; E="-0.23%"; F="11-1-2023"; G="0" },
    @{ Row=14; D="0.001563"; E="0.26%"; F="11-1-2023"; G="0" },
    @{ Row=15; D="0.0006409"; E="1.52%"; F="11-1-2023"; G="0" },
    @{ Row=16; D="0.006025"; E="-0.31%"; F="11-1-2023"; G="0" },
    @{ Row=17; D="3.465"; E="0.26%"; F="11-1-2023"; G="0" },
    @{ Row=18; E="-1.13%"; F="11-1-2023"; G="0" },
    @{ Row=19; E="0.03%"; F="11-1-2023"; G="0" },
    @{ Row=20; F="11-1-2023"; G="0" },
    @{ Row=21; D="0.1337"; E="1.69%"; F="11-1-2023"; G="0" },
    @{ Row=22; D="3.911"; E="0.23%"; F="11-1-2023"; G="0" },
    @{ Row=23; D="0.04413"; E="1.06%"; F="11-1-2023"; G="0" },
    @{ Row=24; D="0.001177"; E="0.10%"; F="11-1-2023"; G="0" },
    @{ Row=25; E="4.97%"; F="11-1-2023"; G="0" },
    @{ Row=26; E="0.11%"; F="11-1-2023"; G="0" },
    @{ Row=27; D="0.0001937"; E="14.20%"; F="11-1-2023"; G="0" },
    @{ Row=28; F="11-1-2023"; G="0" },
    @{ Row=29; F="11-1-2023"; G="0" },
    @{ Row=30; F="11-1-2023"; G="0" },
    @{ Row=31; F="11-1-2023"; G="0" },
    @{ Row=32; F="11-1-2023"; G="0" },
    @{ Row=33; F="11-1-2023"; G="0" },
    @{ Row=34; F="11-1-2023"; G="0" },
    @{ Row=35; F="11-1-2023"; G="0" },
    @{ Row=36; F="11-1-2023"; G="0" },
    @{ Row=37; F="11-1-2023"; G="0" },
    @{ Row=38; F="11-1-2023"; G="0" },
    @{ Row=39; F="11-1-2023"; G="0" },
    @{ Row=40; D="0.04151"; E="2.76%"; F="11-1-2023"; G="0" },
    @{ Row=41; D="0.006810"; E="2.93%"; F="11-1-2023"; G="0" },
    @{ Row=42; D="0.1176"; E="1.05%"; F="11-1-2023"; G="0" },
    @{ Row=43; E="7.85%"; F="11-1-2023"; G="0" },
    @{ Row=44; D="0.01163"; E="-7.85%"; F="11-1-2023"; G="0" },
    @{ Row=45; D="0.00005155"; E="-3.01%"; F="11-1-2023"; G="0" },
    @{ Row=46; E="-36.89%"; F="11-1-2023"; G="0" },
    @{ Row=47; D="0.02025"; E="-22.07%"; F="11-1-2023"; G="0" },
    @{ Row=48; F="11-1-2023"; G="0" },
    @{ Row=49; F="11-1-2023"; G="0" },
    @{ Row=50; F="11-1-2023"; G="0" },
    @{ Row=51; F="11-1-2023"; G="0" },
)

foreach ($row in $updates) {
    foreach ($col in @("D", "E", "F", "G")) {
        if ($row.ContainsKey($col)) {
            $cell = $ws.Cells.Item($row.Row, $colMap[$col])
            $cell.NumberFormat = "@"
            $cell.Value = $row[$col]
        }
    }
}
